$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "57.492.47"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "3.087.49"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "516.07"
$ws.Range("E5").Value = "  +0.44%  "
Set-TextValue "D6" "141.10"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "3.612.28"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  -5.08%  "
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "57.574.80"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "3.087.57"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("E18").Value = "  -0.38%  "
Set-TextValue "D19" "13.18"
$ws.Range("E19").Value = "  -1.45%  "
Set-TextValue "D20" "8.15"
$ws.Range("E20").Value = "  +0.44%  "
Set-TextValue "D21" "334.40"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -1.02%  "
Set-TextValue "D24" "65.90"
$ws.Range("E24").Value = "  +0.79%  "
Set-TextValue "D25" "0.170"
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "0.0₃0911"
$ws.Range("E27").Value = "  +3.05%  "
Set-TextValue "D28" "6.39"
$ws.Range("E28").Value = "  -4.86%  "
Set-TextValue "D29" "7.12"
$ws.Range("E29").Value = "  -1.04%  "
Set-TextValue "D30" "1.82"
$ws.Range("E30").Value = "  +0.53%  "
Set-TextValue "D31" "20.84"
Set-TextValue "D32" "1.16"
$ws.Range("E32").Value = "  -3.38%  "
Set-TextValue "D33" "154.88"
$ws.Range("E33").Value = "  +1.56%  "
Set-TextValue "D34" "28.02"
$ws.Range("E34").Value = "  +11.27%  "
Set-TextValue "D35" "4.55"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("E36").Value = "  -1.03%  "
Set-TextValue "D37" "1.27"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D39").Value = "3.128.31"
$ws.Range("E39").Value = "  +1.70%  "
Set-TextValue "D40" "36.79"
$ws.Range("E40").Value = "  +0.22%  "
Set-TextValue "D41" "3.87"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "2.287.82"
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("E46").Value = "  -0.96%  "
Set-TextValue "D47" "0.941"
$ws.Range("E47").Value = "  -0.74%  "
Set-TextValue "D48" "20.04"
$ws.Range("E48").Value = "  -0.11%  "
Set-TextValue "D49" "5.89"
$ws.Range("E49").Value = "  -3.45%  "
Set-TextValue "D50" "254.05"
$ws.Range("E50").Value = "  +7.50%  "
Set-TextValue "D51" "0.0876"
$ws.Range("E51").Value = "  +1.22%  "
